# Refresh the cryptos price/volume snapshot (GitHub Actions style update).
# Column D ("Price") values are text that often look numeric (e.g.
# "40.152.61", "1.00", "14.20"); a leading apostrophe forces Excel to
# keep them as literal text instead of silently coercing to a Number
# and dropping significant trailing/grouping characters.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range('D2').Value = '''40.152.61'
$ws.Range('E2').Value = '  +0.99%  '
$ws.Range('D3').Value = '''2.235.31'
$ws.Range('E3').Value = '  -0.33%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''293.42'
$ws.Range('E5').Value = '  -1.69%  '
$ws.Range('D6').Value = '''88.47'
$ws.Range('E6').Value = '  +5.44%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('E9').Value = '  -0.11%  '
$ws.Range('D10').Value = '''31.08'
$ws.Range('E10').Value = '  +1.59%  '
$ws.Range('D11').Value = '''0.0791'
$ws.Range('E11').Value = '  +1.18%  '
$ws.Range('D12').Value = '''47.69'
$ws.Range('E12').Value = '  +0.65%  '
$ws.Range('D14').Value = '''6.47'
$ws.Range('E14').Value = '  +1.83%  '
$ws.Range('D15').Value = '''2.576.01'
$ws.Range('E15').Value = '  -0.37%  '
$ws.Range('D16').Value = '''14.20'
$ws.Range('E16').Value = '  -0.67%  '
$ws.Range('D17').Value = '''2.243.86'
$ws.Range('E17').Value = '  +0.33%  '
$ws.Range('D18').Value = '''0.737'
$ws.Range('E18').Value = '  +1.77%  '
$ws.Range('D19').Value = '''40.093.25'
$ws.Range('E19').Value = '  +0.94%  '
$ws.Range('D20').Value = '''11.62'
$ws.Range('E20').Value = '  +10.41%  '
$ws.Range('E21').Value = '  +0.90%  '
$ws.Range('E22').Value = '  +0.89%  '
$ws.Range('D23').Value = '''66.15'
$ws.Range('E23').Value = '  +1.20%  '
$ws.Range('D24').Value = '''236.76'
$ws.Range('E24').Value = '  +3.25%  '
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('D26').Value = '''2.49'
$ws.Range('E26').Value = '  +2.20%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').Value = '''23.00'
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('E29').Value = '  +0.81%  '
$ws.Range('D30').Value = '''9.33'
$ws.Range('E30').Value = '  +1.35%  '
$ws.Range('D31').Value = '''33.05'
$ws.Range('E31').Value = '  +0.88%  '
$ws.Range('D32').Value = '''153.02'
$ws.Range('E32').Value = '  +1.90%  '
$ws.Range('D33').Value = '''0.999'
$ws.Range('E33').Value = '  -0.15%  '
$ws.Range('D34').Value = '''4.99'
$ws.Range('E34').Value = '  +2.11%  '
$ws.Range('E35').Value = '  +2.35%  '
$ws.Range('E36').Value = '  -2.13%  '
$ws.Range('D37').Value = '''2.86'
$ws.Range('E37').Value = '  +6.61%  '
$ws.Range('D38').Value = '''16.23'
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('D39').Value = '''0.113'
$ws.Range('E39').Value = '  +0.97%  '
$ws.Range('E40').Value = '  +2.74%  '
$ws.Range('E41').Value = '  +2.61%  '
$ws.Range('D42').Value = '''2.106.15'
$ws.Range('E42').Value = '  +9.26%  '
$ws.Range('D43').Value = '''3.85'
$ws.Range('E43').Value = '  +2.62%  '
$ws.Range('D44').Value = '''2.19'
$ws.Range('E44').Value = '  +6.93%  '
$ws.Range('E45').Value = '  +2.26%  '
$ws.Range('D46').Value = '''10.11'
$ws.Range('E46').Value = '  +10.39%  '
$ws.Range('D47').Value = '''18.09'
$ws.Range('E47').Value = '  +7.83%  '
$ws.Range('E48').Value = '  +1.80%  '
$ws.Range('D49').Value = '''2.444.85'
$ws.Range('E49').Value = '  -0.39%  '
$ws.Range('D50').Value = '''71.47'
$ws.Range('E50').Value = '  -0.29%  '
$ws.Range('E51').Value = '  +5.72%  '
